# Apply "Update with Correct Forecast output" changes to the
# "Forecast Comparison" worksheet:
#   - insert a new "Week_Start_Date" column between Week and ASIN
#   - strip the leading zero from single-digit week labels (W01 -> W1, ...)
#   - populate the new column with the week's start date (as text)
#   - convert the is_holiday_week column values to real booleans

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B (ASIN), shifting everything right.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Week_Start_Date"

# Make sure the new date column is stored as plain text so values like
# "2025-01-05" are not auto-converted into date serial numbers.
$ws.Range("B2:B17").NumberFormat = "@"

$weeks = @("W1","W2","W3","W4","W5","W6","W7","W8","W9","W10","W11","W12","W13","W14","W15","W16")
$weekStartDates = @("2025-01-05","2025-01-12","2025-01-19","2025-01-26","2025-02-02","2025-02-09","2025-02-16","2025-02-23","2025-03-02","2025-03-09","2025-03-16","2025-03-23","2025-03-30","2025-04-06","2025-04-13","2025-04-20")

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2
    $ws.Range("A$row").Value = $weeks[$i]
    $ws.Range("B$row").Value = $weekStartDates[$i]
}

# The former column I (is_holiday_week) is now column J after the insert.
# Re-write its values as real booleans instead of 0/1 numbers.
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J$row").Value = $false
}

$wb.Save()
